$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organization")

$ws.Range("A3").Value = "Google"
$ws.Range("A4").Value = "Wipro"

$ws.Range("A5").Select()
